$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add initial level 2 dialog/tutorial rows.
# Row 31 is filled first, then row 30, to match the shared-string insertion order.
$ws.Range("A31").Value = "level_2_intro_1"
$ws.Range("B31").Value = "Hi"

$ws.Range("A30").Value = "card_drag_instruct"
$ws.Range("B30").Value = "Drag a card to an empty slot."

$ws.Range("B31").Select()
